# Apply updated dSF (column F) values as part of a data re-pull / recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -5
    4  = -2
    5  = -2
    6  = -5
    7  = -2
    9  = 3
    10 = -7
    13 = -6
    14 = 5
    15 = -4
    17 = 2
    19 = -2
    21 = -6
    22 = -3
    23 = -2
    24 = -2
    25 = -4
    26 = 2
    27 = 1
    28 = -2
    31 = -5
    36 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
